$d = $word.ActiveDocument

# =====================================================================
# Change 1: "Version 3.4.1" -> "Version 3.4.2"
#   Original runs: ... <w:r><w:t>4</w:t></w:r><w:r><w:t>.1</w:t></w:r>
#   Target runs:   ... <w:r><w:t>4</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t>2</w:t></w:r>
# =====================================================================
$verFind = $d.Content
$verFound = $verFind.Find.Execute(".1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($verFound) {
    $verRange = $d.Range($verFind.Start, $verFind.End)
    $verRange.Delete()
    $dotIns = $d.Range($verRange.Start, $verRange.Start)
    $dotIns.InsertAfter(".")
    $twoIns = $d.Range($dotIns.End, $dotIns.End)
    $twoIns.InsertAfter("2")
}

# =====================================================================
# Change 2: comparator formula correction  13 -> 31
#   "= year * 372 + month * 13 + day"
#     -> "= year * 372 + month * 31 + day (version 3.4.2 correction 13 -> 31)"
#   with "31" and the parenthetical note colored red, matching highlighting
#   kept on the plain-text runs.
# =====================================================================
$formulaPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "year \* 372") {
        $formulaPara = $p
        break
    }
}

$thirteen = $formulaPara.Range
$thirteen.Find.Execute("13", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$thirteen.Text = "31"
$thirteen.Font.Color = 255

$tailAnchor = $formulaPara.Range
$tailInsert = $d.Range($tailAnchor.End - 1, $tailAnchor.End - 1)
$tailInsert.InsertAfter(" (version 3.4.2 correction 13 -> 31)")

$highlightRange = $formulaPara.Range
$highlightRange.Find.Execute(" (version 3.4.2 correction 13 -> 31)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$highlightRange.HighlightColorIndex = 8

$noteRange = $formulaPara.Range
$noteRange.Find.Execute("(version 3.4.2 correction 13 -> 31)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$noteRange.Font.Color = 255

# ---------------------------------------------------------------------
# Move the _GoBack bookmark so it now sits at the end of this paragraph
# (right after the note, before the paragraph mark) instead of at the
# end of the next paragraph.
#
# NOTE: adding a collapsed bookmark exactly at (paragraph.End - 1) can
# mis-place it, so a one-character placeholder is appended first, the
# bookmark is anchored just before the placeholder (a safe, non-edge
# position), and then the placeholder is removed -- the bookmark stays
# put, now correctly sitting at the true end of the paragraph content.
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$placeholderAnchor = $formulaPara.Range
$placeholderAt = $d.Range($placeholderAnchor.End - 1, $placeholderAnchor.End - 1)
$placeholderAt.InsertAfter("Z")

$afterPlaceholder = $formulaPara.Range
$bookmarkAt = $afterPlaceholder.End - 2
$bookmarkRange = $d.Range($bookmarkAt, $bookmarkAt)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$placeholderChar = $d.Range($bookmarkAt, $bookmarkAt + 1)
$placeholderChar.Delete()

# =====================================================================
# Change 3: merge the trailing runs of the next paragraph into a single
# run now that the bookmark no longer separates them.
#   "...min_date" + " " + "and all other conditions must be met." + bookmark
#     -> "...min_date" + " and all other conditions must be met."
# =====================================================================
$conditionsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "and all other conditions") {
        $conditionsPara = $p
        break
    }
}

$tailText = $conditionsPara.Range
$tailText.Find.Execute(" and all other conditions must be met.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tailText.Delete()
$mergedIns = $d.Range($tailText.Start, $tailText.Start)
$mergedIns.InsertAfter(" and all other conditions must be met.")
$mergedIns.HighlightColorIndex = 8

# =====================================================================
# Change 4: header page-number field result  7 -> 8
# =====================================================================
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers(1)
    if ($hdr.Exists) {
        $hdrRange = $hdr.Range
        if ($hdrRange.Find.Execute("7", $true, $false, $false, $false, $false, $true, 1, $false, "8", 2)) {
            # replaced
        }
    }
}
